# LOM3118.docx edit.
#
# The course-description texts that used to sit under "Objetivos",
# "Programa resumido", "Programa", the three "Avaliacao" answers (Metodo /
# Criterio / Norma de recuperacao) and "Bibliografia" get redistributed:
#
#   * The text that was under "Objetivos" moves into the docente-roster
#     paragraph (prepended before the roster, as its first line).
#   * The text that was under "Programa" moves into the docente-roster
#     paragraph too (second prepended line).
#   * The three "Avaliacao" answers and the "Bibliografia" text likewise
#     move into the docente-roster paragraph, each becoming its own
#     prepended line (in that order).
#   * "Objetivos" itself receives the text that used to be under
#     "Programa resumido".
#   * The roster's final six names (Luiz, Maria, Miguel, Sandra,
#     Sebastiao, Sergio) are trimmed off the roster paragraph and
#     redistributed one-each into "Programa resumido", "Programa", the
#     three "Avaliacao" answers, and "Bibliografia" -- the slots the
#     relocated texts vacated.
#
# All source text is read back live from the document (instead of being
# retyped as literals here) so accents / punctuation round-trip exactly.

$d = $word.ActiveDocument
$lb = [char]11   # Word's soft line break char; serializes as <w:br/>
$cr = [char]13   # paragraph-mark char Range.Text reports at a paragraph end

function Get-ParaText($index) {
    return $d.Paragraphs.Item($index).Range.Text.TrimEnd($cr).TrimEnd($lb)
}

# ---------------------------------------------------------------------
# 0) Snapshot all the texts that are about to move, before anything is
#    edited (paragraph indices below are the ORIGINAL, 1-based layout).
# ---------------------------------------------------------------------
$objetivosOld        = Get-ParaText 6     # "Objetivos" answer
$programaResumidoOld = Get-ParaText 10    # "Programa resumido" answer
$programaOld         = Get-ParaText 12    # "Programa" answer
$bibliografiaOld     = Get-ParaText 16    # "Bibliografia" answer

$p14 = $d.Paragraphs.Item(14)             # "Avaliacao" answer paragraph

$rMetodo = $p14.Range.Duplicate
$rMetodo.Find.Execute("Método: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$metodoLabelEnd = $rMetodo.End

$rCriterio = $p14.Range.Duplicate
$rCriterio.Find.Execute("Critério: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$criterioLabelStart = $rCriterio.Start
$criterioLabelEnd = $rCriterio.End

$rNorma = $p14.Range.Duplicate
$rNorma.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$normaLabelStart = $rNorma.Start
$normaLabelEnd = $rNorma.End

$p14End = $p14.Range.End

$metodoOld   = $d.Range($metodoLabelEnd, $criterioLabelStart).Text.TrimEnd($cr).TrimEnd($lb)
$criterioOld = $d.Range($criterioLabelEnd, $normaLabelStart).Text.TrimEnd($cr).TrimEnd($lb)
$normaOld    = $d.Range($normaLabelEnd, $p14End).Text.TrimEnd($cr).TrimEnd($lb)

# ---------------------------------------------------------------------
# 1) Docente-roster paragraph: prepend the six relocated texts (each as
#    its own run ending in a line break) before the existing roster.
#    Insert in reverse order at the fixed paragraph-start offset so the
#    runs end up in forward (reading) order.
# ---------------------------------------------------------------------
$listPara = $d.Paragraphs.Item(8)
$insertPos = $listPara.Range.Start

$blocks = @($objetivosOld, $programaOld, $metodoOld, $criterioOld, $normaOld, $bibliografiaOld)
for ($i = $blocks.Length - 1; $i -ge 0; $i--) {
    $ins = $d.Range($insertPos, $insertPos)
    $ins.InsertBefore($blocks[$i] + $lb)
}

# ---------------------------------------------------------------------
# 2) Trim the roster: remove the final six names (and the line break
#    right before the first of them), leaving Hugo as the paragraph's
#    last run (no trailing break).
# ---------------------------------------------------------------------
$hugoRange = $listPara.Range.Duplicate
$hugoRange.Find.Execute("984972 - Hugo Ricardo Zschommler Sandim", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$delStart = $hugoRange.End
$delEnd = $d.Paragraphs.Item(8).Range.End
$d.Range($delStart, $delEnd).Delete() | Out-Null

# ---------------------------------------------------------------------
# 3) "Objetivos" answer now receives the "Programa resumido" text.
# ---------------------------------------------------------------------
$d.Paragraphs.Item(6).Range.Text = $programaResumidoOld

# ---------------------------------------------------------------------
# 4) The four vacated single-answer paragraphs ("Programa resumido",
#    "Programa", "Bibliografia") and the two remaining "Avaliacao"
#    answer runs receive the relocated roster names.
# ---------------------------------------------------------------------
$d.Paragraphs.Item(10).Range.Text = "1176388 - Luiz Tadeu Fernandes Eleno"
$d.Paragraphs.Item(12).Range.Text = "7459752 - Maria Ismenia Sodero Toledo Faria"

$p14b = $d.Paragraphs.Item(14)
$rMetodo2 = $p14b.Range.Duplicate
$rMetodo2.Find.Execute("Método: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$metodoLabelEnd2 = $rMetodo2.End

$rCriterio2 = $p14b.Range.Duplicate
$rCriterio2.Find.Execute("Critério: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$criterioLabelStart2 = $rCriterio2.Start
$criterioLabelEnd2 = $rCriterio2.End

$rNorma2 = $p14b.Range.Duplicate
$rNorma2.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$normaLabelStart2 = $rNorma2.Start
$normaLabelEnd2 = $rNorma2.End

$p14End2 = $p14b.Range.End

$d.Range($metodoLabelEnd2, $criterioLabelStart2).Text = "5840622 - Miguel Justino Ribeiro Barboza" + $lb

# Recompute offsets: the Método answer's length changed, so re-find the
# "Critério: " / "Norma de recuperação: " labels after that edit.
$rCriterio3 = $p14b.Range.Duplicate
$rCriterio3.Find.Execute("Critério: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$criterioLabelEnd3 = $rCriterio3.End

$rNorma3 = $p14b.Range.Duplicate
$rNorma3.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$normaLabelStart3 = $rNorma3.Start
$normaLabelEnd3 = $rNorma3.End

$p14End3 = $p14b.Range.End

$d.Range($criterioLabelEnd3, $normaLabelStart3).Text = "2166002 - Sandra Giacomin Schneider" + $lb

$rNorma4 = $p14b.Range.Duplicate
$rNorma4.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$normaLabelEnd4 = $rNorma4.End
$p14End4 = $p14b.Range.End

$d.Range($normaLabelEnd4, $p14End4 - 1).Text = "1922320 - Sebastiao Ribeiro"

$d.Paragraphs.Item(16).Range.Text = "5840793 - Sérgio Schneider"
